$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Ford Tourneo Custom"
$ws.Range("C1").Value = "ZEEKR 001"
$ws.Range("D1").Value = "MAXUS MIFA 7"
$ws.Range("E1").Value = "VW Passat"
$ws.Range("F1").Value = "Škoda Kodiaq"
$ws.Range("G1").Value = "BMW X2"
$ws.Range("H1").Value = "Renault Rafale HEV"
$ws.Range("I1").Value = "Mercedes-Benz E-Class"
$ws.Range("J1").Value = "Suzuki Swift"
$ws.Range("K1").Value = "Dacia Duster"
$ws.Range("L1").Value = "Renault Espace"
$ws.Range("M1").Value = "Toyota C-HR"
$ws.Range("N1").Value = "Honda CR-V"
$ws.Range("O1").Value = "NIO EL6"
